$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (YZA567): register a maintenance entry
$ws.Range("L3").Value = 12
$ws.Range("P3").Value = "10/10/2024"
$ws.Range("Q3").Value = 12

# Row 7 (KLM789): truck departs on a trip, fuel type corrected
$ws.Range("D7").Value = "En viaje"
$ws.Range("E7").Value = "Disel"
$ws.Range("H7").Value = "21/01/1900"
